$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pasos")
$ws.Activate()

# Insert a new row at position 4 (between the old row 3 "buscar ciudad" step and
# the old row 4 "presiono boton buscar" step), pushing everything else down.
$ws.Rows("4").Insert()

# Copy formatting (borders/fonts/number formats) from row 3 into the freshly
# inserted row 4 so it matches the surrounding table look.
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Update step text (column B) and expected-result text (column C) ---
# Order matches the order these strings were first introduced by the author
# (controls the resulting shared-string table ordering).
$ws.Range("C3").Value = "Se muestra un listado de ciudades que coinciden con <Ciudad1>"

$ws.Range("B4").Value = "Seleccionar del listado la ciudad <Ciudad1>"

$ws.Range("B2").Value = "Ingresar a la pagina principal."

$ws.Range("B3").Value = "Ingresar <Ciudad1> en el campo nombre de ciudad"

$ws.Range("B5").Value = "Presionar el boton buscar"
$ws.Range("C5").Value = "Se carga la pagina BusquedaPlayas, con todas las playas de que pertenecen a la ciudad <Ciudad1> disponibles en un mapa. Se carga la informacion de las playas en la grilla de playas."

$ws.Range("B6").Value = "Selecciono <TipoPlaya1> en el campo Tipo de Playa"
$ws.Range("B7").Value = "Selecciono <TipoVehiculo1> en el campo Tipo de Vehiculo"
$ws.Range("B8").Value = "Ingreso <Precio1> en el campo precio desde"
$ws.Range("B9").Value = "Ingreso <PrecioHasta> en el campo precio hasta"
$ws.Range("B10").Value = "Ingreso <Horario1> en el campo hora desde"
$ws.Range("B11").Value = "Ingreso <Horario2> en el campo hora hasta"
$ws.Range("B12").Value = "Selecciono <DiasDeAtencion1> en el campo Dias de atencion"

$ws.Range("C13").Value = "Se muestran en el mapa solo las playas que cumplen con los filtros seleccionados]"

# Update the selection to match the author's final cursor position.
$ws.Range("B2:C13").Select()

$wb.Save()
